$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.670.56"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.291.81"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'95.68"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("D6").Value = "'267.94"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  -4.41%  "
$ws.Range("D10").Value = "'44.82"
$ws.Range("E10").Value = "  -7.54%  "
$ws.Range("D11").Value = "'0.0938"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  -4.73%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "2.632.22"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'15.18"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "'0.853"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "2.291.42"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "43.558.68"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'72.30"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "'2.45"
$ws.Range("E22").Value = "  +4.67%  "
$ws.Range("D23").Value = "'235.31"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'9.07"
$ws.Range("E24").Value = "  -14.55%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "'11.22"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").Value = "'40.42"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").Value = "'175.00"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "'22.01"
$ws.Range("E32").Value = "  +3.45%  "
$ws.Range("D33").Value = "'0.0884"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0357"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").Value = "'4.39"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'3.31"
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("E40").Value = "  -6.56%  "
$ws.Range("E41").Value = "  +6.03%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'12.11"
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.34"
$ws.Range("E43").Value = "  +14.87%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'64.32"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "'8.81"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "'98.25"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'1.50"
$ws.Range("E50").Value = "  +5.25%  "
$ws.Range("D51").Value = "2.511.56"
$ws.Range("E51").Value = "  +1.22%  "
